$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab_5a_Indikatoren")

$ws.Range("E2").Value = "a) Material and social deprivation"
$ws.Range("L2").Value = "a) Materielle und soziale Deprivation"
$ws.Range("M2").Value = "a) Material and social deprivation"
$ws.Range("E3").Value = "b) Severe material and social deprivation"
$ws.Range("L3").Value = "b) Erhebliche materielle und soziale Deprivation"
$ws.Range("M3").Value = "b) Severe material and social deprivation"
$ws.Range("L8").Value = "c) Jugendliche"
$ws.Range("M8").Value = "c) Adolescents"
$ws.Range("L9").Value = "d) Erwachsene"
$ws.Range("M9").Value = "d) Adults"
$ws.Range("L18").Value = "a) 0- bis 2-Jährige"
$ws.Range("M18").Value = "a) 0 to 2-year-olds"
$ws.Range("L19").Value = "b) 3- bis 5-Jährige"
$ws.Range("M19").Value = "b) 3 to 5-year-olds"
$ws.Range("L23").Value = "b) In Aufsichtsräten der börsennotierten und paritätisch mitbestimmten Unternehmen"
$ws.Range("M23").Value = "b) On supervisory boards of listed and fully co-determined companies"
$ws.Range("L24").Value = "c) Im öffentlichen Dienst des Bundes"
$ws.Range("M24").Value = "c) In management positions in the federal civil service"
$ws.Range("L25").Value = "Kinder, deren Väter Elterngeld bezogen haben"
$ws.Range("M25").Value = "XXXKinder, deren Väter Elterngeld bezogen haben"
$ws.Range("L31").Value = "a) Mit einer neuen oder verbesserten Trinkwasserversorgung"
$ws.Range("M31").Value = "XXXPeople gaining first-time or upgraded access to drinking water"
$ws.Range("L32").Value = "b) Mit einer neuen oder verbesserten Basissanitärversorgung oder Abwasserbehandlung"
$ws.Range("L33").Value = "Mit einem neuen oder hochwertigeren Zugang zur Trinkwasserversorgung oder Anschluss zur Sanitärversorgung"
$ws.Range("M33").Value = "XXXPeople gaining first-time or upgraded access to drinking water or sanitation"
$ws.Range("L34").Value = "a) Endenergieproduktivität"
$ws.Range("M34").Value = "a) Final energy productivity"
$ws.Range("L35").Value = "b) Primärenergieverbrauch"
$ws.Range("M35").Value = "b) Primary energy consumption"
$ws.Range("L39").Value = "a) Staatsdefizit"
$ws.Range("M39").Value = "a) Government deficit"
$ws.Range("L40").Value = "b) Strukturelles Defizit"
$ws.Range("M40").Value = "b) Structural deficit"
$ws.Range("L44").Value = "a) 20- bis 64-Jährige"
$ws.Range("M44").Value = "a) 20 to 64-year-olds"
$ws.Range("L45").Value = "b) 55- bis 64-Jährige"
$ws.Range("M45").Value = "b) 55 to 64-year-olds"
$ws.Range("L60").Value = "ba) Direkter und indirekter Rohstoffeinsatz"
$ws.Range("M60").Value = "ba) Direct and indirect use of raw materials"
$ws.Range("L61").Value = "bb) Direkter und indirekter Energieverbrauch"
$ws.Range("M61").Value = "bb) Direct and indirect energy consumption"
$ws.Range("L62").Value = "bc) Direkte und indirekte CO2-Emissionen"
$ws.Range("M62").Value = "bc) Direct and indirect CO2 emissions"
$ws.Range("L65").Value = "a) Nachhaltige Vergabeverfahren"
$ws.Range("L66").Value = "b) CO₂-Emissionen je Fahrleistungen der Kraftfahrzeuge"
$ws.Range("L67").Value = "c) Nachhaltige Textilbeschaffung"
$ws.Range("L71").Value = "aa) Stickstoffeintrag über die Zuflüsse in die Ostsee"
$ws.Range("L72").Value = "ab) Stickstoffeintrag über die Zuflüsse in die Nordsee"
$ws.Range("L77").Value = "Bilaterale Beiträge der deutschen internationalen Kooperation zum Schutz, nachhaltiger Nutzung und Wiederherstellung von Land (inkl. Wald)"
$ws.Range("L79").Value = "a) Corruption Perception Index in Deutschland"
$ws.Range("M79").Value = "XXXCorruption Perception Index in Germany"
$ws.Range("L80").Value = "b) Corruption Perception Index in den Partnerländern der deutschen Entwicklungszusammenarbeit"
$ws.Range("M80").Value = "XXXNumber of partner countries for German development cooperation with improved CPI scores compared with 2012"
